$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new PO Forecast sheet and populate it BEFORE moving it ---
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Range("A2").Value = 44948.99999999999
$wsForecast.Range("B2").Value = 49
$wsForecast.Range("C2").Value = -75.59530120665319
$wsForecast.Range("D2").Value = 166.3745075216089
$wsForecast.Range("A3").Value = 44962.99999999999
$wsForecast.Range("B3").Value = 48
$wsForecast.Range("C3").Value = -66.49878592273816
$wsForecast.Range("D3").Value = 172.8534961407292
$wsForecast.Range("A4").Value = 44969.99999999999
$wsForecast.Range("B4").Value = 48
$wsForecast.Range("C4").Value = -71.77238817381686
$wsForecast.Range("D4").Value = 178.9360086199443
$wsForecast.Range("A5").Value = 44976.99999999999
$wsForecast.Range("B5").Value = 47
$wsForecast.Range("C5").Value = -73.79119508081078
$wsForecast.Range("D5").Value = 166.538923623922
$wsForecast.Range("A6").Value = 44983.99999999999
$wsForecast.Range("B6").Value = 47
$wsForecast.Range("C6").Value = -77.95448128072319
$wsForecast.Range("D6").Value = 168.9600853998828
$wsForecast.Range("A7").Value = 44990.99999999999
$wsForecast.Range("B7").Value = 46
$wsForecast.Range("C7").Value = -82.13476183018388
$wsForecast.Range("D7").Value = 171.0513259166142
$wsForecast.Range("A8").Value = 44997.99999999999
$wsForecast.Range("B8").Value = 46
$wsForecast.Range("C8").Value = -71.79774103166034
$wsForecast.Range("D8").Value = 165.6731226795948
$wsForecast.Range("A9").Value = 45004.99999999999
$wsForecast.Range("B9").Value = 45
$wsForecast.Range("C9").Value = -77.02904464173373
$wsForecast.Range("D9").Value = 175.8527979332153
$wsForecast.Range("A10").Value = 45011.99999999999
$wsForecast.Range("B10").Value = 45
$wsForecast.Range("C10").Value = -81.46662350784891
$wsForecast.Range("D10").Value = 162.6816329407339
$wsForecast.Range("A11").Value = 45025.99999999999
$wsForecast.Range("B11").Value = 44
$wsForecast.Range("C11").Value = -80.86349339113323
$wsForecast.Range("D11").Value = 170.9849346769948
$wsForecast.Range("A12").Value = 45032.99999999999
$wsForecast.Range("B12").Value = 43
$wsForecast.Range("C12").Value = -83.58096166070897
$wsForecast.Range("D12").Value = 167.7233704020926
$wsForecast.Range("A13").Value = 45039.99999999999
$wsForecast.Range("B13").Value = 43
$wsForecast.Range("C13").Value = -74.88996182609456
$wsForecast.Range("D13").Value = 175.8437557905773
$wsForecast.Range("A14").Value = 45046.99999999999
$wsForecast.Range("B14").Value = 42
$wsForecast.Range("C14").Value = -74.35587407005364
$wsForecast.Range("D14").Value = 161.4743550962844
$wsForecast.Range("A15").Value = 45053.99999999999
$wsForecast.Range("B15").Value = 42
$wsForecast.Range("C15").Value = -83.20211215949404
$wsForecast.Range("D15").Value = 157.6303851548278
$wsForecast.Range("A16").Value = 45067.99999999999
$wsForecast.Range("B16").Value = 41
$wsForecast.Range("C16").Value = -86.30052238114594
$wsForecast.Range("D16").Value = 159.7921063678874
$wsForecast.Range("A17").Value = 45109.99999999999
$wsForecast.Range("B17").Value = 38
$wsForecast.Range("C17").Value = -78.82735104725477
$wsForecast.Range("D17").Value = 159.2924790329773
$wsForecast.Range("A18").Value = 45158.99999999999
$wsForecast.Range("B18").Value = 35
$wsForecast.Range("C18").Value = -86.91387355372724
$wsForecast.Range("D18").Value = 152.1183750144732
$wsForecast.Range("A19").Value = 45165.99999999999
$wsForecast.Range("B19").Value = 34
$wsForecast.Range("C19").Value = -92.3091934932407
$wsForecast.Range("D19").Value = 154.8455632887962
$wsForecast.Range("A20").Value = 45172.99999999999
$wsForecast.Range("B20").Value = 34
$wsForecast.Range("C20").Value = -84.66418399018181
$wsForecast.Range("D20").Value = 159.959723506
$wsForecast.Range("A21").Value = 45179.99999999999
$wsForecast.Range("B21").Value = 33
$wsForecast.Range("C21").Value = -85.45606979243864
$wsForecast.Range("D21").Value = 158.9897697061442
$wsForecast.Range("A22").Value = 45186.99999999999
$wsForecast.Range("B22").Value = 33
$wsForecast.Range("C22").Value = -95.96503831842567
$wsForecast.Range("D22").Value = 149.8078326294235
$wsForecast.Range("A23").Value = 45193.99999999999
$wsForecast.Range("B23").Value = 32
$wsForecast.Range("C23").Value = -90.00503979726062
$wsForecast.Range("D23").Value = 156.3420143835893
$wsForecast.Range("A24").Value = 45207.99999999999
$wsForecast.Range("B24").Value = 31
$wsForecast.Range("C24").Value = -88.156618594058
$wsForecast.Range("D24").Value = 145.8658426010983
$wsForecast.Range("A25").Value = 45214.99999999999
$wsForecast.Range("B25").Value = 31
$wsForecast.Range("C25").Value = -91.41004945378451
$wsForecast.Range("D25").Value = 157.7232413288338
$wsForecast.Range("A26").Value = 45235.99999999999
$wsForecast.Range("B26").Value = 29
$wsForecast.Range("C26").Value = -83.47622895198386
$wsForecast.Range("D26").Value = 167.3385325345693
$wsForecast.Range("A27").Value = 45242.99999999999
$wsForecast.Range("B27").Value = 29
$wsForecast.Range("C27").Value = -89.99593857273359
$wsForecast.Range("D27").Value = 151.1410985793694
$wsForecast.Range("A28").Value = 45249.99999999999
$wsForecast.Range("B28").Value = 28
$wsForecast.Range("C28").Value = -93.77531444638237
$wsForecast.Range("D28").Value = 148.0994361668444
$wsForecast.Range("A29").Value = 45256.99999999999
$wsForecast.Range("B29").Value = 28
$wsForecast.Range("C29").Value = -98.98377222474163
$wsForecast.Range("D29").Value = 153.4322269957472
$wsForecast.Range("A30").Value = 45263.99999999999
$wsForecast.Range("B30").Value = 27
$wsForecast.Range("C30").Value = -98.30666578761306
$wsForecast.Range("D30").Value = 149.2270823270797
$wsForecast.Range("A31").Value = 45270.99999999999
$wsForecast.Range("B31").Value = 27
$wsForecast.Range("C31").Value = -96.63811023270442
$wsForecast.Range("D31").Value = 142.1115214019671
$wsForecast.Range("A32").Value = 45277.99999999999
$wsForecast.Range("B32").Value = 26
$wsForecast.Range("C32").Value = -95.00542172335263
$wsForecast.Range("D32").Value = 147.3793000010791
$wsForecast.Range("A33").Value = 45298.99999999999
$wsForecast.Range("B33").Value = 25
$wsForecast.Range("C33").Value = -96.26351867884298
$wsForecast.Range("D33").Value = 148.8082044536791
$wsForecast.Range("A34").Value = 45305.99999999999
$wsForecast.Range("B34").Value = 24
$wsForecast.Range("C34").Value = -101.1513999841982
$wsForecast.Range("D34").Value = 146.4044313248765
$wsForecast.Range("A35").Value = 45312.99999999999
$wsForecast.Range("B35").Value = 24
$wsForecast.Range("C35").Value = -98.04219878427679
$wsForecast.Range("D35").Value = 140.0607568223085
$wsForecast.Range("A36").Value = 45340.99999999999
$wsForecast.Range("B36").Value = 22
$wsForecast.Range("C36").Value = -89.93633935977097
$wsForecast.Range("D36").Value = 140.6346469486815
$wsForecast.Range("A37").Value = 45354.99999999999
$wsForecast.Range("B37").Value = 21
$wsForecast.Range("C37").Value = -98.44840539267355
$wsForecast.Range("D37").Value = 153.9040932491547
$wsForecast.Range("A38").Value = 45361.99999999999
$wsForecast.Range("B38").Value = 20
$wsForecast.Range("C38").Value = -103.0331998972459
$wsForecast.Range("D38").Value = 137.9762782832216
$wsForecast.Range("A39").Value = 45368.99999999999
$wsForecast.Range("B39").Value = 20
$wsForecast.Range("C39").Value = -104.4137915578647
$wsForecast.Range("D39").Value = 144.9590178519088
$wsForecast.Range("A40").Value = 45375.99999999999
$wsForecast.Range("B40").Value = 20
$wsForecast.Range("C40").Value = -105.058045675779
$wsForecast.Range("D40").Value = 139.848250543592
$wsForecast.Range("A41").Value = 45389.99999999999
$wsForecast.Range("B41").Value = 19
$wsForecast.Range("C41").Value = -93.28162851182489
$wsForecast.Range("D41").Value = 147.1903907650027
$wsForecast.Range("A42").Value = 45396.99999999999
$wsForecast.Range("B42").Value = 18
$wsForecast.Range("C42").Value = -101.2030496800437
$wsForecast.Range("D42").Value = 139.6911472520035
$wsForecast.Range("A43").Value = 45515.99999999999
$wsForecast.Range("B43").Value = 10
$wsForecast.Range("C43").Value = -113.7362455885783
$wsForecast.Range("D43").Value = 134.9222971313502
$wsForecast.Range("A44").Value = 45564.99999999999
$wsForecast.Range("B44").Value = 6
$wsForecast.Range("C44").Value = -115.9092129085031
$wsForecast.Range("D44").Value = 124.1511660062127
$wsForecast.Range("A45").Value = 45571.99999999999
$wsForecast.Range("B45").Value = 6
$wsForecast.Range("C45").Value = -121.237699209932
$wsForecast.Range("D45").Value = 128.9457865281843
$wsForecast.Range("A46").Value = 45578.99999999999
$wsForecast.Range("B46").Value = 5
$wsForecast.Range("C46").Value = -113.7440688275308
$wsForecast.Range("D46").Value = 129.0736535661317
$wsForecast.Range("A47").Value = 45585.99999999999
$wsForecast.Range("B47").Value = 5
$wsForecast.Range("C47").Value = -114.6711600812744
$wsForecast.Range("D47").Value = 127.7645356653808
$wsForecast.Range("A48").Value = 45592.99999999999
$wsForecast.Range("B48").Value = 4
$wsForecast.Range("C48").Value = -118.194613346872
$wsForecast.Range("D48").Value = 126.239381243136
$wsForecast.Range("A49").Value = 45599.99999999999
$wsForecast.Range("B49").Value = 4
$wsForecast.Range("C49").Value = -118.0500267156799
$wsForecast.Range("D49").Value = 123.4717255597387
$wsForecast.Range("A50").Value = 45606.99999999999
$wsForecast.Range("B50").Value = 4
$wsForecast.Range("C50").Value = -111.6326622898972
$wsForecast.Range("D50").Value = 127.023062782834
$wsForecast.Range("A51").Value = 45613.99999999999
$wsForecast.Range("B51").Value = 3
$wsForecast.Range("C51").Value = -124.1605901547125
$wsForecast.Range("D51").Value = 126.6148779509915
$wsForecast.Range("A52").Value = 45620.99999999999
$wsForecast.Range("B52").Value = 3
$wsForecast.Range("C52").Value = -111.6354212577613
$wsForecast.Range("D52").Value = 128.0498331021545

# Now move the populated sheet to the end (last tab position)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Move($null, $lastSheet)

Write-Output "Done"